$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cells (outside the used A1:D32 range) used to stage text values.
# Excel's Range.Value setter auto-coerces comma-grouped numeric-looking
# strings (and date-like strings such as "Sep 2022") into numbers/dates.
# To keep these as literal text (matching the original shared-string
# cells, with no style/number-format changes), we append a throwaway
# letter so the initial entry is unambiguously text, trim it off with a
# LEFT()/LEN() text formula (which evaluates to a text result), then
# copy that computed text and paste-special "values only" into the
# real target cells. This avoids creating any new cell styles.
function Set-TextValue {
    param([string[]]$targets, [string]$text)

    $ws.Range("Y1").Value = $text + "X"
    $ws.Range("Z1").Formula = "=LEFT(Y1,LEN(Y1)-1)"
    $ws.Range("Z1").Copy()
    foreach ($t in $targets) {
        $ws.Range($t).PasteSpecial(-4163)
    }
    $ws.Range("Y1:Z1").Clear()
}

# Header: "Aug 2022" -> "Sep 2022"
Set-TextValue @("D3") "Sep 2022"

# ASET LANCAR detail rows
Set-TextValue @("D6") "7,198,061"
Set-TextValue @("D7") "670,013,184"
Set-TextValue @("D8") "4,692,951,888"
Set-TextValue @("D9") "-1,520,457,188"

# JUMLAH ASET LANCAR / JUMLAH ASET (same figure, two cells)
Set-TextValue @("D10", "D19") "3,849,705,945"

# Aset Neto Tidak Terikat / JUMLAH ASET NETO (same figure, two cells)
Set-TextValue @("D30", "D31") "3,712,330,587"

# JUMLAH LIABILITAS DAN ASET NETO (now balances with JUMLAH ASET)
Set-TextValue @("D32") "3,849,705,945"
